$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 18 de Mayo de 2020 a las 10:35"

# Polonia / Ucrania swapped places (row 34 <-> row 35), data follows the row
$ws.Range("A34").Value = "Polonia"
$ws.Range("B34").Value = 18746
$ws.Range("C34").Value = 217
$ws.Range("D34").Value = 7451
$ws.Range("E34").Value = 10366
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = 4
$ws.Range("H34").Value = 929

$ws.Range("A35").Value = "Ucrania"
$ws.Range("B35").Value = 18616
$ws.Range("C35").Value = 325
$ws.Range("D35").Value = 5276
$ws.Range("E35").Value = 12805
$ws.Range("F35").Value = 0
$ws.Range("G35").Value = 21
$ws.Range("H35").Value = 535

# Israel (row 38)
$ws.Range("B38").Value = 16621
$ws.Range("C38").Value = 4
$ws.Range("D38").Value = 13014
$ws.Range("E38").Value = 3335

# Filipinas (row 44)
$ws.Range("B44").Value = 12718
$ws.Range("C44").Value = 205
$ws.Range("D44").Value = 2729
$ws.Range("E44").Value = 9158
$ws.Range("G44").Value = 7
$ws.Range("H44").Value = 831

# Dinamarca (row 48)
$ws.Range("B48").Value = 10968
$ws.Range("C48").Value = 41
$ws.Range("E48").Value = 1194

# Lituania (row 89)
$ws.Range("B89").Value = 1547
$ws.Range("C89").Value = 6
$ws.Range("E89").Value = 491
$ws.Range("G89").Value = 3
$ws.Range("H89").Value = 59

# Taiwan (row 132)
$ws.Range("D132").Value = 398
$ws.Range("E132").Value = 35
